$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Record 3 pages written on 2021-09-28 (row 18, column I)
$ws.Range("I18").Value = 3

# Recalculate all dependent formulas (B5, B11, D11, B12, I31, etc.)
$excel.Calculate()

# Selection ended up on I19 (the cell below, e.g. after pressing Enter)
$ws.Range("I19").Select()
